$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 705867065
$ws.Range("B2").Value = 79999999999
$ws.Range("C2").Value = 2858164260
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = "Кораблев Кирил Олегович"
